$d = $word.ActiveDocument

# The last paragraph of the document ends with "...lished without taking
# more than 3 trips. " followed by the _GoBack bookmark. We need to:
#   1) end that paragraph right after "...3 trips. " (drop the bookmark
#      from it)
#   2) add a blank paragraph
#   3) add a bold+underlined "Problem 2:" paragraph
#   4) add a blank bold+underlined paragraph
#   5) add a new paragraph that starts with a tab and states the
#      constraints of problem 2, carrying the _GoBack bookmark at its end
$last = $d.Paragraphs.Last

$newBody = '<w:p w14:paraId="38CB47AC" w14:textId="3F48DD65" w:rsidR="00700D7A" w:rsidRDefault="00700D7A" w:rsidP="00F447CC"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">There is only one solution to this problem, which is to take the parrot first, return for the seed and then pick the parrot back up, return for the cat while leaving the parrot, and then to return for the parrot. </w:t></w:r><w:r w:rsidR="00365580"><w:t>My first plan was to take the seed first, since it did not interact with anything, but then I remembered that the cat and the parrot could not be left behind. I looked into taking each one by itself, but soon realized that the journey could not be accomp</w:t></w:r><w:r w:rsidR="003F0C7C"><w:t xml:space="preserve">lished without taking more than 3 trips. </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Problem 2:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">The constraints of the problem are that there are 20 socks in the drawer of various colors in different amounts. The socks can only be selected in the dark, so determining their color can only be done after selection. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>{0}<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@ -f $newBody

$last.Range.InsertXML($xml)
